# Update countries & provincias Spain
# - Reorders three country rows (keeping label/data pairs correctly attached)
# - Refreshes a handful of per-country case statistics
# - Bumps the "last updated" timestamp in the footer

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Footer timestamp: 22:05 -> 22:35
# ---------------------------------------------------------------------------
$footer = $ws.Cells.Find("Datos actualizados")
$footer.Value = $footer.Value().Replace("22:05", "22:35")

# ---------------------------------------------------------------------------
# 2) Tayikistan moves up, right after Mali / before Kenia.
#    Kenia, Crucero, Uruguay each shift down one row; Tayikistan's row gets
#    freshly updated figures.
# ---------------------------------------------------------------------------
$keniaRow = $ws.Cells.Find("Kenia").Row

$rKenia    = $keniaRow
$rCrucero  = $keniaRow + 1
$rUruguay  = $keniaRow + 2
$rTayik    = $keniaRow + 3

$keniaVals   = @($ws.Cells.Item($rKenia,2).Value(),   $ws.Cells.Item($rKenia,3).Value(),   $ws.Cells.Item($rKenia,4).Value(),   $ws.Cells.Item($rKenia,5).Value(),   $ws.Cells.Item($rKenia,6).Value(),   $ws.Cells.Item($rKenia,7).Value(),   $ws.Cells.Item($rKenia,8).Value())
$cruceroVals = @($ws.Cells.Item($rCrucero,2).Value(), $ws.Cells.Item($rCrucero,3).Value(), $ws.Cells.Item($rCrucero,4).Value(), $ws.Cells.Item($rCrucero,5).Value(), $ws.Cells.Item($rCrucero,6).Value(), $ws.Cells.Item($rCrucero,7).Value(), $ws.Cells.Item($rCrucero,8).Value())
$uruguayVals = @($ws.Cells.Item($rUruguay,2).Value(), $ws.Cells.Item($rUruguay,3).Value(), $ws.Cells.Item($rUruguay,4).Value(), $ws.Cells.Item($rUruguay,5).Value(), $ws.Cells.Item($rUruguay,6).Value(), $ws.Cells.Item($rUruguay,7).Value(), $ws.Cells.Item($rUruguay,8).Value())

# Row that used to be Kenia now becomes Tayikistan, with brand-new totals.
$ws.Cells.Item($rKenia,1).Value = "Tayikistan"
$ws.Cells.Item($rKenia,2).Value = 729
$ws.Cells.Item($rKenia,3).Value = 68
$ws.Cells.Item($rKenia,4).Value = 0
$ws.Cells.Item($rKenia,5).Value = 708
$ws.Cells.Item($rKenia,6).Value = 0
$ws.Cells.Item($rKenia,7).Value = 0
$ws.Cells.Item($rKenia,8).Value = 21

# Row that used to be Crucero now becomes Kenia, carrying Kenia's old data.
$ws.Cells.Item($rCrucero,1).Value = "Kenia"
$ws.Cells.Item($rCrucero,2).Value = $keniaVals[0]
$ws.Cells.Item($rCrucero,3).Value = $keniaVals[1]
$ws.Cells.Item($rCrucero,4).Value = $keniaVals[2]
$ws.Cells.Item($rCrucero,5).Value = $keniaVals[3]
$ws.Cells.Item($rCrucero,6).Value = $keniaVals[4]
$ws.Cells.Item($rCrucero,7).Value = $keniaVals[5]
$ws.Cells.Item($rCrucero,8).Value = $keniaVals[6]

# Row that used to be Uruguay now becomes Crucero, carrying Crucero's old data.
$ws.Cells.Item($rUruguay,1).Value = "Crucero"
$ws.Cells.Item($rUruguay,2).Value = $cruceroVals[0]
$ws.Cells.Item($rUruguay,3).Value = $cruceroVals[1]
$ws.Cells.Item($rUruguay,4).Value = $cruceroVals[2]
$ws.Cells.Item($rUruguay,5).Value = $cruceroVals[3]
$ws.Cells.Item($rUruguay,6).Value = $cruceroVals[4]
$ws.Cells.Item($rUruguay,7).Value = $cruceroVals[5]
$ws.Cells.Item($rUruguay,8).Value = $cruceroVals[6]

# Row that used to be Tayikistan now becomes Uruguay, carrying Uruguay's old data.
$ws.Cells.Item($rTayik,1).Value = "Uruguay"
$ws.Cells.Item($rTayik,2).Value = $uruguayVals[0]
$ws.Cells.Item($rTayik,3).Value = $uruguayVals[1]
$ws.Cells.Item($rTayik,4).Value = $uruguayVals[2]
$ws.Cells.Item($rTayik,5).Value = $uruguayVals[3]
$ws.Cells.Item($rTayik,6).Value = $uruguayVals[4]
$ws.Cells.Item($rTayik,7).Value = $uruguayVals[5]
$ws.Cells.Item($rTayik,8).Value = $uruguayVals[6]

# ---------------------------------------------------------------------------
# 3) Nueva Caledonia moves up, right before Belice (their "activos"/"muertes"
#    figures swap along with them).
# ---------------------------------------------------------------------------
$rBelice = $ws.Cells.Find("Belice").Row
$rNuevaCaledonia = $rBelice + 1

$beliceD = $ws.Cells.Item($rBelice,4).Value()
$beliceH = $ws.Cells.Item($rBelice,8).Value()
$ncD = $ws.Cells.Item($rNuevaCaledonia,4).Value()
$ncH = $ws.Cells.Item($rNuevaCaledonia,8).Value()

$ws.Cells.Item($rBelice,1).Value = "Nueva Caledonia"
$ws.Cells.Item($rBelice,4).Value = $ncD
$ws.Cells.Item($rBelice,8).Value = $ncH

$ws.Cells.Item($rNuevaCaledonia,1).Value = "Belice"
$ws.Cells.Item($rNuevaCaledonia,4).Value = $beliceD
$ws.Cells.Item($rNuevaCaledonia,8).Value = $beliceH

# ---------------------------------------------------------------------------
# 4) San Bartolome moves up, right before Sahara Occidental (figures for both
#    are identical, so this is purely a label/ordering swap).
# ---------------------------------------------------------------------------
$rSahara = $ws.Cells.Find("Sahara Occidental").Row
$rSanBartolome = $rSahara + 1

$ws.Cells.Item($rSahara,1).Value = "San Bartolome"
$ws.Cells.Item($rSanBartolome,1).Value = "Sahara Occidental"

# ---------------------------------------------------------------------------
# 5) Refreshed case counts for a handful of countries.
# ---------------------------------------------------------------------------
$rUSA = $ws.Cells.Find("Estados Unidos").Row
$ws.Cells.Item($rUSA,2).Value = 1402913
$ws.Cells.Item($rUSA,3).Value = 17079
$ws.Cells.Item($rUSA,4).Value = 276383
$ws.Cells.Item($rUSA,5).Value = 1043371
$ws.Cells.Item($rUSA,7).Value = 1364
$ws.Cells.Item($rUSA,8).Value = 83159

$rBrasil = $ws.Cells.Find("Brasil").Row
$ws.Cells.Item($rBrasil,2).Value = 173141
$ws.Cells.Item($rBrasil,3).Value = 3998
$ws.Cells.Item($rBrasil,5).Value = 93692
$ws.Cells.Item($rBrasil,7).Value = 440
$ws.Cells.Item($rBrasil,8).Value = 12065

$rAlemania = $ws.Cells.Find("Alemania").Row
$ws.Cells.Item($rAlemania,2).Value = 173034
$ws.Cells.Item($rAlemania,3).Value = 458
$ws.Cells.Item($rAlemania,5).Value = 18116
$ws.Cells.Item($rAlemania,7).Value = 57
$ws.Cells.Item($rAlemania,8).Value = 7718

$rYibuti = $ws.Cells.Find("Republica de Yibuti").Row
$ws.Cells.Item($rYibuti,2).Value = 1256
$ws.Cells.Item($rYibuti,3).Value = 29
$ws.Cells.Item($rYibuti,4).Value = 886
$ws.Cells.Item($rYibuti,5).Value = 367

$rRuanda = $ws.Cells.Find("Ruanda").Row
$ws.Cells.Item($rRuanda,2).Value = 286
$ws.Cells.Item($rRuanda,3).Value = 1
$ws.Cells.Item($rRuanda,4).Value = 153
$ws.Cells.Item($rRuanda,5).Value = 133
